$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Casos de Uso")

$ws1.Range("F7").Value = "Hecho"
$ws1.Range("F8").Value = "Hecho"
$ws1.Range("G8").Value = 1
$ws1.Range("K8").Value = 1

$ws1.Range("G9").Value = 1
$ws1.Range("K9").Value = 2

$ws1.Range("G10").Value = 1
$ws1.Range("K10").Value = 2

$ws1.Range("G11").Value = 1
$ws1.Range("K11").Value = 1
